$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while preserving it as literal text (no numeric
# auto-conversion) and without leaving behind any NumberFormat/style change.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "35.794.08"
Set-TextValue "E2" "  +3.84%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.871.49"
Set-TextValue "E3" "  +3.22%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.35%  "

# Row 5 - BNB
Set-TextValue "D5" "232.83"
Set-TextValue "E5" "  +3.22%  "

# Row 6 - XRP
Set-TextValue "D6" "0.616"
Set-TextValue "E6" "  +3.61%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.35%  "

# Row 8 - Solana
Set-TextValue "D8" "42.71"
Set-TextValue "E8" "  +11.50%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.313"
Set-TextValue "E9" "  +7.76%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +3.67%  "

# Row 11 - TRON
Set-TextValue "E11" "  +4.31%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.144.80"
Set-TextValue "E12" "  +3.35%  "

# Row 13 - Chainlink
Set-TextValue "D13" "11.75"
Set-TextValue "E13" "  +4.73%  "

# Row 14 - was Polygon, now WrappedEther
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.876.13"
Set-TextValue "E14" "  +3.55%  "

# Row 15 - was WrappedEther, now Polygon
Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.688"
Set-TextValue "E15" "  +8.66%  "

# Row 16 - Polkadot
Set-TextValue "E16" "  +8.43%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "35.818.27"
Set-TextValue "E17" "  +3.91%  "

# Row 18 - Litecoin
Set-TextValue "D18" "70.87"
Set-TextValue "E18" "  +3.82%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0₃0810"
Set-TextValue "E19" "  +4.71%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "249.67"
Set-TextValue "E20" "  +2.71%  "

# Row 21 - Avalanche
Set-TextValue "D21" "12.54"
Set-TextValue "E21" "  +11.59%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.82"
Set-TextValue "E22" "  +16.67%  "

# Row 23 - Dai
Set-TextValue "E23" "  +0.37%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  +1.90%  "

# Row 25 - Monero
Set-TextValue "D25" "172.04"
Set-TextValue "E25" "  +1.08%  "

# Row 26 - Cosmos
Set-TextValue "D26" "8.12"
Set-TextValue "E26" "  +3.88%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "18.06"
Set-TextValue "E27" "  +2.94%  "

# Row 28 - Stellar
Set-TextValue "E28" "  +2.25%  "

# Row 29 - PancakeSwap
Set-TextValue "D29" "1.44"
Set-TextValue "E29" "  +17.19%  "

# Row 30 - BinanceUSD
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.33%  "

# Row 31 - EURNeutrino
Set-TextValue "D31" "3.302.96"
Set-TextValue "E31" "  +35.94%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.0557"
Set-TextValue "E32" "  +7.61%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.99"
Set-TextValue "E33" "  +5.13%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "E34" "  +6.55%  "

# Row 35 - LidoDAOToken
Set-TextValue "E35" "  +5.31%  "

# Row 36 - Aave
Set-TextValue "D36" "99.50"
Set-TextValue "E36" "  +21.72%  "

# Row 37 - ImmutableX
Set-TextValue "E37" "  +7.24%  "

# Row 38 - RenderToken
Set-TextValue "E38" "  +7.15%  "

# Row 39 - Maker
Set-TextValue "D39" "1.370.32"
Set-TextValue "E39" "  +0.67%  "

# Row 40 - TrustWalletToken
Set-TextValue "D40" "1.10"
Set-TextValue "E40" "  +3.21%  "

# Row 41 - VeChain
Set-TextValue "E41" "  +6.22%  "

# Row 42 - ARBITRUM
Set-TextValue "E42" "  +8.38%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "15.14"
Set-TextValue "E43" "  +9.68%  "

# Row 44 - WEMIXToken
Set-TextValue "E44" "  +3.99%  "

# Row 45 - HuobiToken
Set-TextValue "E45" "  +1.81%  "

# Row 46 - MXToken
Set-TextValue "D46" "2.83"
Set-TextValue "E46" "  +1.08%  "

# Row 47 - FraxShare
Set-TextValue "D47" "6.34"
Set-TextValue "E47" "  +9.97%  "

# Row 48 - Kaspa
Set-TextValue "D48" "0.0520"
Set-TextValue "E48" "  +2.29%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.042.70"
Set-TextValue "E49" "  +3.38%  "

# Row 50 - Quant
Set-TextValue "D50" "105.53"
Set-TextValue "E50" "  +3.42%  "

# Row 51 - PaxDollar
Set-TextValue "E51" "  +0.40%  "
